# Weekly update: insert a new price record as the first data row for this
# market/category block (row 10 of data -> sheet row 410), pushing the
# existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 410; Excel copies the formatting
# of the row above (keeps the date style on column D) and shifts every
# row below it (410-432) down to (411-433), which also bumps the used
# range to A1:R433.
$ws.Rows.Item(410).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A410").Value = 11
$ws.Range("B410").Value = "Vega Monumental Concepción"
$ws.Range("C410").Value = "Bíobío"
$ws.Range("D410").Value = (Get-Date -Year 2023 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E410").Value = 8
$ws.Range("F410").Value = 100114001
$ws.Range("G410").Value = "Papa"
$ws.Range("H410").Value = "Asterix"
$ws.Range("I410").Value = "1a (cosecha lavada)"
$ws.Range("J410").Value = 4000
$ws.Range("K410").Value = 11500
$ws.Range("L410").Value = 12000
$ws.Range("M410").Value = 11750
$ws.Range("N410").Value = "$/malla 25 kilos"
$ws.Range("O410").Value = "Región de La Araucanía"
$ws.Range("P410").Value = 470
$ws.Range("Q410").Value = 25
$ws.Range("R410").Value = "Hortaliza"
